# OB API for Saturn nit
# Edits the big body textbox ("TextBox 1") on slide 1:
#  - paragraph 1: merge ". However, due to expensive TTP "/"certification "/"schemes "
#    into one run, move "schemes " in front of "as ", and turn the straight
#    quotes around "Wallet" into curly quotes (split into separate runs).
#  - paragraph 3: turn the (unbalanced) straight/curly quote around "Wallet"
#    into a properly closed curly-quoted "Wallet", and italicize "dedicated for ".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 1")
$tr = $sh.TextFrame.TextRange

$ldq = [char]0x201C                   # “
$rdq = [char]0x201D                   # ”

# Work right-to-left (highest character offset first) so offsets to the
# left of each edit stay valid for subsequent operations.

# --- Paragraph 3 ("Saturn is an open, light-weight scheme (including ...") ---

# "dedicated for " -> italic
$tr.Characters(346, 14).Font.Italic = $true

# "Wallet), " -> "Wallet" + right-curly-quote + "), "
$tr.Characters(337, 9).Text = "Wallet" + $rdq + "), "

# --- Paragraph 1 ("In theory Open Banking APIs can support ...") ---

# closing quote: '"' -> right curly quote
$tr.Characters(175, 1).Text = $rdq

# '("Wallet' -> '(' + left-curly-quote + 'Wallet'
$tr.Characters(167, 8).Text = "(" + $ldq + "Wallet"

# "as " -> "schemes as "
$tr.Characters(112, 3).Text = "schemes as "

# ". However, due to expensive TTP " + "certification " + "schemes " -> merged run
$tr.Characters(58, 54).Text = ". However, due to expensive TTP certification "
